$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 454, shifting existing rows 454:547 down to 455:548.
$ws.Rows.Item(454).Insert(-4121)

# Populate the newly inserted row 454 with the new price-report record.
$ws.Range("A454").Value = 4
$ws.Range("B454").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C454").Value = "Los Lagos"
$ws.Range("D454").Value = 44932
$ws.Range("E454").Value = 10
$ws.Range("F454").Value = "Fruta"
$ws.Range("G454").Value = 100103
$ws.Range("H454").Value = "Frutos de hueso (carozo)"
$ws.Range("I454").Value = 100103006
$ws.Range("J454").Value = "Nectarín"
$ws.Range("K454").Value = "Super Queen"
$ws.Range("L454").Value = "Primera"
$ws.Range("M454").Value = 500
$ws.Range("N454").Value = 18000
$ws.Range("O454").Value = 18000
$ws.Range("P454").Value = 18000
$ws.Range("Q454").Value = "`$/caja 14 kilos empedrada"
$ws.Range("R454").Value = "Región de O'Higgins"
$ws.Range("S454").Value = 1286
$ws.Range("T454").Value = 14
